$d = $word.ActiveDocument

# Namespace-wrapped pkg:package header/footer used for Range.InsertXML payloads.
$xmlHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParagraphByText($doc, $text) {
    $want = $text + [char]13
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq $want) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphRuns($doc, $oldText, $innerXml) {
    $p = Get-ParagraphByText $doc $oldText
    if ($p -eq $null) {
        throw "Paragraph with text '$oldText' not found"
    }
    $payload = $xmlHeader + '<w:p>' + $innerXml + '</w:p>' + $xmlFooter
    [void]$p.Range.InsertXML($payload)
}

$langRpr = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

# 1) "Category 1B Description" -> "Description" / " of" / " Group" / " 1" / "B"
$inner1 = ''
$inner1 += '<w:r>' + $langRpr + '<w:t>Description</w:t></w:r>'
$inner1 += '<w:r>' + $langRpr + '<w:t xml:space="preserve"> of</w:t></w:r>'
$inner1 += '<w:r>' + $langRpr + '<w:t xml:space="preserve"> Group</w:t></w:r>'
$inner1 += '<w:r>' + $langRpr + '<w:t xml:space="preserve"> 1</w:t></w:r>'
$inner1 += '<w:r>' + $langRpr + '<w:t>B</w:t></w:r>'
Set-ParagraphRuns $d "Category 1B Description" $inner1

# 2) "Category 1C Description" -> "Description" / " of Group" / " 1C"
$inner2 = ''
$inner2 += '<w:r>' + $langRpr + '<w:t>Description</w:t></w:r>'
$inner2 += '<w:r>' + $langRpr + '<w:t xml:space="preserve"> of Group</w:t></w:r>'
$inner2 += '<w:r>' + $langRpr + '<w:t xml:space="preserve"> 1C</w:t></w:r>'
Set-ParagraphRuns $d "Category 1C Description" $inner2

# 3) "Description of Group 5A" -> single run, drop proofErr wrapping, add lang rPr to pPr & run
$inner3 = '<w:pPr>' + $langRpr + '</w:pPr>'
$inner3 += '<w:r>' + $langRpr + '<w:t>Description of Group 5A</w:t></w:r>'
Set-ParagraphRuns $d "Description of Group 5A" $inner3

# 4) "Description of Group 6A" -> single run, drop proofErr wrapping, add lang rPr to pPr & run
$inner4 = '<w:pPr>' + $langRpr + '</w:pPr>'
$inner4 += '<w:r>' + $langRpr + '<w:t>Description of Group 6A</w:t></w:r>'
Set-ParagraphRuns $d "Description of Group 6A" $inner4

Write-Output "Done"
